$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "1623"
$ws.Range("E4").Value = "14167367.43"

$ws.Range("C6").Value = "1125"
$ws.Range("E6").Value = "8635487.99"

$ws.Range("C15").Value = "307"
$ws.Range("E15").Value = "1176266.84"

$ws.Range("C17").Value = "772"
$ws.Range("E17").Value = "8014452.50"

$ws.Range("E18").Value = "1945770.62"

$ws.Range("C81").Value = "1334"
$ws.Range("E81").Value = "10996644.05"

$ws.Range("C82").Value = "698"
$ws.Range("E82").Value = "5110678.51"
